$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("icf_a.185")
$ws.Activate()

# Update row 12 values from BW12 through CG12 to the new value 1540
$ws.Range("BW12:CG12").Value = 1540

# Update the selected cell/range to reflect CG13 as the new active selection
$ws.Range("CG13").Select()
